$wb = $excel.ActiveWorkbook

$oldGuid = "4ef547bc-ef3a-4dc1-bc6d-0856992128d6"
$newGuid = "4987c566-b142-4352-bddb-92d8c3dc69ee"
$oldHash = "35a9babb923a2d0304a6cc69c79c7b50d51299de"
$newHash = "c4809855ef7909218a54bfdcf61b514fc1587d33"

# Hyperlink target URL is unchanged by this edit (only the cell/display text
# changes) - it still points at the commit that contains the *old* guid.
$linkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e8b4f2f4fca6b26e582c3b9620c4a97f2b5e7bc/e2e/$oldGuid.md"

function Update-Hyperlink($ws, $cellRef, $displayText) {
    $r = $ws.Range($cellRef)
    $r.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($r, $linkTarget, "", "", $displayText) | Out-Null
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
Update-Hyperlink $wsOverview "B2" "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-19 21:01:18"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
Update-Hyperlink $wsZh "A2" "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-19 21:01:14"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
Update-Hyperlink $wsDe "A2" "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
# de-de!H2 ("Latest Handoff Datetime") originally shared the exact same
# string as Overview!G2 ("2016-08-19 21:00:53"); that shared string's text
# is what changed, so both cells move to the new timestamp together.
$wsDe.Range("H2").Value = "2016-08-19 21:01:18"
